# Auto update Excel log
# Appends new sensor-log rows to the PIR, Humidity and Temperature sheets,
# matching the source system's export for 2026-01-28.
#
# Columns (all sheets): A=Date, B=Timestamp, C=Hour, D=Location, E=Value, F=Status
#
# NOTE: several "Value" strings look like other data types to Excel's
# auto-detection (an ISO date "2026-01-28", or a percentage like "86.7%").
# Force those cells to Text format *before* assigning so they are stored
# as literal strings instead of being coerced into a date serial /
# percentage number.

$wb = $excel.ActiveWorkbook

function Add-LogRows {
    param($ws, $rows)

    foreach ($row in $rows) {
        $r = $row[0]

        $dateCell = $ws.Cells.Item($r, 1)
        $dateCell.NumberFormat = "@"
        $dateCell.Value = $row[1]

        $ws.Cells.Item($r, 2).Value = $row[2]
        $ws.Cells.Item($r, 3).Value = $row[3]
        $ws.Cells.Item($r, 4).Value = $row[4]

        $valueCell = $ws.Cells.Item($r, 5)
        $valueCell.NumberFormat = "@"
        $valueCell.Value = $row[5]

        $ws.Cells.Item($r, 6).Value = $row[6]
    }
}

# ---------------------------------------------------------------------
# PIR sheet: rows 321-333
# ---------------------------------------------------------------------
$pirWs = $wb.Worksheets.Item("PIR")
$pirRows = @(
    @(321, "2026-01-28", "12:31:13", "12:00", "Bathroom", "No Motion", "Inactive"),
    @(322, "2026-01-28", "12:31:14", "12:00", "Bathroom", "No Motion", "Inactive"),
    @(323, "2026-01-28", "12:31:19", "12:00", "Bathroom", "No Motion", "Inactive"),
    @(324, "2026-01-28", "12:31:26", "12:00", "Bathroom", "No Motion", "Inactive"),
    @(325, "2026-01-28", "12:31:30", "12:00", "Bathroom", "No Motion", "Inactive"),
    @(326, "2026-01-28", "12:31:34", "12:00", "Bathroom", "No Motion", "Inactive"),
    @(327, "2026-01-28", "12:31:39", "12:00", "Bathroom", "No Motion", "Inactive"),
    @(328, "2026-01-28", "12:31:46", "12:00", "Bathroom", "No Motion", "Inactive"),
    @(329, "2026-01-28", "12:31:50", "12:00", "Bathroom", "No Motion", "Inactive"),
    @(330, "2026-01-28", "12:31:54", "12:00", "Bathroom", "No Motion", "Inactive"),
    @(331, "2026-01-28", "12:31:59", "12:00", "Bathroom", "No Motion", "Inactive"),
    @(332, "2026-01-28", "12:32:04", "12:00", "Bathroom", "No Motion", "Inactive"),
    @(333, "2026-01-28", "12:32:10", "12:00", "Bathroom", "No Motion", "Inactive")
)
Add-LogRows $pirWs $pirRows

# ---------------------------------------------------------------------
# Humidity sheet: rows 299-310
# ---------------------------------------------------------------------
$humidityWs = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
    @(299, "2026-01-28", "12:31:11", "12:00", "Bathroom", "86.7%", "Active"),
    @(300, "2026-01-28", "12:31:16", "12:00", "Bathroom", "86.7%", "Active"),
    @(301, "2026-01-28", "12:31:20", "12:00", "Bathroom", "87.5%", "Active"),
    @(302, "2026-01-28", "12:31:24", "12:00", "Bathroom", "87.5%", "Active"),
    @(303, "2026-01-28", "12:31:28", "12:00", "Bathroom", "87.5%", "Active"),
    @(304, "2026-01-28", "12:31:32", "12:00", "Bathroom", "87.5%", "Active"),
    @(305, "2026-01-28", "12:31:36", "12:00", "Bathroom", "86.6%", "Active"),
    @(306, "2026-01-28", "12:31:40", "12:00", "Bathroom", "87.5%", "Active"),
    @(307, "2026-01-28", "12:31:44", "12:00", "Bathroom", "87.4%", "Active"),
    @(308, "2026-01-28", "12:31:52", "12:00", "Bathroom", "87.3%", "Active"),
    @(309, "2026-01-28", "12:31:56", "12:00", "Bathroom", "86.4%", "Active"),
    @(310, "2026-01-28", "12:32:08", "12:00", "Bathroom", "86.3%", "Active")
)
Add-LogRows $humidityWs $humidityRows

# ---------------------------------------------------------------------
# Temperature sheet: rows 299-310
# ---------------------------------------------------------------------
$temperatureWs = $wb.Worksheets.Item("Temperature")
$temperatureRows = @(
    @(299, "2026-01-28", "12:31:12", "12:00", "Bathroom", "22.9C", "Active"),
    @(300, "2026-01-28", "12:31:17", "12:00", "Bathroom", "23.0C", "Active"),
    @(301, "2026-01-28", "12:31:21", "12:00", "Bathroom", "23.0C", "Active"),
    @(302, "2026-01-28", "12:31:25", "12:00", "Bathroom", "23.0C", "Active"),
    @(303, "2026-01-28", "12:31:29", "12:00", "Bathroom", "23.0C", "Active"),
    @(304, "2026-01-28", "12:31:33", "12:00", "Bathroom", "23.0C", "Active"),
    @(305, "2026-01-28", "12:31:37", "12:00", "Bathroom", "23.0C", "Active"),
    @(306, "2026-01-28", "12:31:41", "12:00", "Bathroom", "23.0C", "Active"),
    @(307, "2026-01-28", "12:31:45", "12:00", "Bathroom", "23.0C", "Active"),
    @(308, "2026-01-28", "12:31:53", "12:00", "Bathroom", "23.0C", "Active"),
    @(309, "2026-01-28", "12:31:57", "12:00", "Bathroom", "23.1C", "Active"),
    @(310, "2026-01-28", "12:32:09", "12:00", "Bathroom", "23.0C", "Active")
)
Add-LogRows $temperatureWs $temperatureRows
